# This workbook tracks daily Apio (celery) prices at Macroferia Regional de
# Talca. A new daily record is inserted at the top of the data block (row 92,
# right after the header row), pushing every existing record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92; this shifts rows 92..227 down to 93..228
# and keeps everything above row 92 untouched.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row with the new price record.
$ws.Range("A92").Value = 5
$ws.Range("B92").Value = "Macroferia Regional de Talca"
$ws.Range("C92").Value = "Maule"
$ws.Range("D92").Value = 44803
$ws.Range("E92").Value = 7
$ws.Range("F92").Value = 100112017
$ws.Range("G92").Value = "Apio"
$ws.Range("H92").Value = "Americana (o)"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 550
$ws.Range("K92").Value = 9000
$ws.Range("L92").Value = 10000
$ws.Range("M92").Value = 9818
$ws.Range("N92").Value = "`$/docena de matas"
$ws.Range("O92").Value = "Provincia del Elquí"
$ws.Range("P92").Value = 1636
$ws.Range("Q92").Value = 6
$ws.Range("R92").Value = "Hortaliza"
